$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 24.16710790014911
$ws.Range("C2").Value = 17.21645978077995
$ws.Range("D2").Value = 6.024638979581598
$ws.Range("E2").Value = 6.558167693945621
$ws.Range("G2").Value = 3.69894102941331
$ws.Range("M2").Value = 19.74809138707936
$ws.Range("N2").Value = 22.3590627182041

$ws.Range("B3").Value = 23.42995828404617
$ws.Range("C3").Value = 16.48961580866731
$ws.Range("D3").Value = 5.915984051653375
$ws.Range("E3").Value = 6.529222824216391
$ws.Range("G3").Value = 3.705528270886033
$ws.Range("M3").Value = 19.41990098967344
$ws.Range("N3").Value = 22.23343510075372

$ws.Range("B4").Value = 22.9751649266292
$ws.Range("C4").Value = 16.03380604724623
$ws.Range("D4").Value = 5.850536758844544
$ws.Range("E4").Value = 6.513224319621575
$ws.Range("G4").Value = 3.709765500730165
$ws.Range("M4").Value = 19.22291003855044
$ws.Range("N4").Value = 22.15707664679928

$ws.Range("B5").Value = 22.78959908996511
$ws.Range("C5").Value = 15.84597833201119
$ws.Range("D5").Value = 5.82422010273438
$ws.Range("E5").Value = 6.507152535549103
$ws.Range("G5").Value = 3.711540927887756
$ws.Range("M5").Value = 19.14387087207875
$ws.Range("N5").Value = 22.12616244835083

$ws.Range("B6").Value = 22.75878140610268
$ws.Range("C6").Value = 15.81467380910337
$ws.Range("D6").Value = 5.819872611553483
$ws.Range("E6").Value = 6.506171400412879
$ws.Range("G6").Value = 3.711838687097172
$ws.Range("M6").Value = 19.13082401048177
$ws.Range("N6").Value = 22.12104159471513

$ws.Range("B7").Value = 22.97266284017457
$ws.Range("C7").Value = 16.03128092512261
$ws.Range("D7").Value = 5.850180365528781
$ws.Range("E7").Value = 6.513140618739086
$ws.Range("G7").Value = 3.709789247116203
$ws.Range("M7").Value = 19.22183895317757
$ws.Range("N7").Value = 22.15665889832235

$ws.Range("B8").Value = 23.91357468018613
$ws.Range("C8").Value = 16.9679981618407
$ws.Range("D8").Value = 5.986932879289034
$ws.Range("E8").Value = 6.547820566037451
$ws.Range("G8").Value = 3.70117249883181
$ws.Range("M8").Value = 19.63405397098619
$ws.Range("N8").Value = 22.31558928402432

$ws.Range("B9").Value = 25.72849709810658
$ws.Range("C9").Value = 18.71730009657852
$ws.Range("D9").Value = 6.263605828064874
$ws.Range("E9").Value = 6.629800076920504
$ws.Range("G9").Value = 3.685790659373864
$ws.Range("M9").Value = 20.47378667878162
$ws.Range("N9").Value = 22.63310062134662

$ws.Range("B10").Value = 27.02766569669353
$ws.Range("C10").Value = 19.93533123039377
$ws.Range("D10").Value = 6.469988507793244
$ws.Range("E10").Value = 6.698375190414309
$ws.Range("G10").Value = 3.675395372604525
$ws.Range("M10").Value = 21.10383406589022
$ws.Range("N10").Value = 22.86950115167006

$ws.Range("B11").Value = 27.60833672216116
$ws.Range("C11").Value = 20.47255521224907
$ws.Range("D11").Value = 6.564125350723801
$ws.Range("E11").Value = 6.731335323069743
$ws.Range("G11").Value = 3.670859007901035
$ws.Range("M11").Value = 21.3920391472484
$ws.Range("N11").Value = 22.97762362517658

$ws.Range("B12").Value = 27.82652360631662
$ws.Range("C12").Value = 20.67340692928316
$ws.Range("D12").Value = 6.599774554334813
$ws.Range("E12").Value = 6.744065319666583
$ws.Range("G12").Value = 3.669168572069185
$ws.Range("M12").Value = 21.50130157910364
$ws.Range("N12").Value = 23.01864125335427

$ws.Range("B13").Value = 27.77961199338788
$ws.Range("C13").Value = 20.63026708057065
$ws.Range("D13").Value = 6.592097324839267
$ws.Range("E13").Value = 6.741312711705929
$ws.Range("G13").Value = 3.669531424350434
$ws.Range("M13").Value = 21.47776599535162
$ws.Range("N13").Value = 23.00980424095037

$ws.Range("B14").Value = 27.62632232872977
$ws.Range("C14").Value = 20.48913199963211
$ws.Range("D14").Value = 6.567058393033451
$ws.Range("E14").Value = 6.732377673848394
$ws.Range("G14").Value = 3.670719387550778
$ws.Range("M14").Value = 21.4010263012989
$ws.Range("N14").Value = 22.98099674421707

$ws.Range("B15").Value = 27.5322004860108
$ws.Range("C15").Value = 20.40234197142329
$ws.Range("D15").Value = 6.551720531458876
$ws.Range("E15").Value = 6.726936935438946
$ws.Range("G15").Value = 3.671450607113902
$ws.Range("M15").Value = 21.35403433745088
$ws.Range("N15").Value = 22.96336066177886

$ws.Range("B16").Value = 26.98949073858909
$ws.Range("C16").Value = 19.89986894169055
$ws.Range("D16").Value = 6.463838526602972
$ws.Range("E16").Value = 6.696256259595799
$ws.Range("G16").Value = 3.675695683237707
$ws.Range("M16").Value = 21.08502230785357
$ws.Range("N16").Value = 22.86244598748441

$ws.Range("B17").Value = 26.65375044918295
$ws.Range("C17").Value = 19.58717823483692
$ws.Range("D17").Value = 6.409965296547476
$ws.Range("E17").Value = 6.67788323998519
$ws.Range("G17").Value = 3.678348995868134
$ws.Range("M17").Value = 20.92032542624214
$ws.Range("N17").Value = 22.80068045838435

$ws.Range("B18").Value = 26.45968093829499
$ws.Range("C18").Value = 19.40574727252541
$ws.Range("D18").Value = 6.379004583307905
$ws.Range("E18").Value = 6.667482028763039
$ws.Range("G18").Value = 3.67989324824508
$ws.Range("M18").Value = 20.82575213109812
$ws.Range("N18").Value = 22.7652100093014

$ws.Range("B19").Value = 26.3938142854957
$ws.Range("C19").Value = 19.34405180055258
$ws.Range("D19").Value = 6.368527318021834
$ws.Range("E19").Value = 6.663989094968178
$ws.Range("G19").Value = 3.680419229561112
$ws.Range("M19").Value = 20.7937612442683
$ws.Range("N19").Value = 22.75321014223381

$ws.Range("B20").Value = 26.68959162343984
$ws.Range("C20").Value = 19.62062939348901
$ws.Range("D20").Value = 6.415697784990328
$ws.Range("E20").Value = 6.67982188630188
$ws.Range("G20").Value = 3.678064671211579
$ws.Range("M20").Value = 20.93784230744882
$ws.Range("N20").Value = 22.80724983571918

$ws.Range("B21").Value = 27.67139498338416
$ws.Range("C21").Value = 20.53065804840229
$ws.Range("D21").Value = 6.574413159323742
$ws.Range("E21").Value = 6.734995404034749
$ws.Range("G21").Value = 3.670369713047738
$ws.Range("M21").Value = 21.42356401499906
$ws.Range("N21").Value = 22.98945627542455

$ws.Range("B22").Value = 28.3030629355256
$ws.Range("C22").Value = 21.11029957313935
$ws.Range("D22").Value = 6.678137274478038
$ws.Range("E22").Value = 6.772501476859548
$ws.Range("G22").Value = 3.665500106752329
$ws.Range("M22").Value = 21.7416967740831
$ws.Range("N22").Value = 23.10896529937555

$ws.Range("B23").Value = 27.96691022087727
$ws.Range("C23").Value = 20.80236322088698
$ws.Range("D23").Value = 6.622789342066025
$ws.Range("E23").Value = 6.752353180431358
$ws.Range("G23").Value = 3.668084611731771
$ws.Range("M23").Value = 21.57187393887892
$ws.Range("N23").Value = 23.04514535063472

$ws.Range("B24").Value = 26.67339107608084
$ws.Range("C24").Value = 19.60551128576492
$ws.Range("D24").Value = 6.413106087800631
$ws.Range("E24").Value = 6.678944920316547
$ws.Range("G24").Value = 3.678193155722847
$ws.Range("M24").Value = 20.92992256515576
$ws.Range("N24").Value = 22.80427969683783

$ws.Range("B25").Value = 25.24246425610249
$ws.Range("C25").Value = 18.25498960292135
$ws.Range("D25").Value = 6.188063895605329
$ws.Range("E25").Value = 6.606142203704371
$ws.Range("G25").Value = 3.689791472991616
$ws.Range("M25").Value = 20.243881766341
$ws.Range("N25").Value = 22.54662313741058
